$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as the third data row (row 3),
# pushing all the existing records (previously rows 3-115) down by one
# row (now rows 4-116). Row 2 (the most recent record before the edit)
# stays untouched.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new "Puerro" price record.
$ws.Range("A3").Value = 9
$ws.Range("B3").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C3").Value = "Metropolitana"
$ws.Range("D3").Value = 44860
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 100112005
$ws.Range("G3").Value = "Puerro"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 9000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 9571
$ws.Range("N3").Value = "`$/paquete 20 unidades"
$ws.Range("O3").Value = "Provincia de Chacabuco"
$ws.Range("P3").Value = 479
$ws.Range("Q3").Value = 20
$ws.Range("R3").Value = "Hortaliza"
